$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59 values for columns A (1) through DB (106)
$values = @(43569,0,0,5,2,0,0,3,109,2,2,3,0,0,8,4,3,9,4,0,2,5,2,4,0,9,9,15,0,0,4,0,0,0,0,9,5,5,2,3,0,95,0,4,7,5,0,5,0,2,0,0,0,0,0,2,0,0,0,2,5,2,4,3,0,0,0,0,2,0,0,0,0,0,0,2,2,0,0,0,0,0,0,2,0,0,0,0,0,0,0,0,0,0,0,0,0,0,2,0,6,0,0,3,2,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(59, $i + 1).Value = $values[$i]
}

# Fix selection/active cell on the sheet view (bug fix @ regions)
$ws.Range("A59").Select()

Write-Host "Row 59 populated and selection fixed"
